# Apply updated "想去人数" (interested count) values to column F
# on sheet "展览" (Exhibitions) and sheet "全部类型" (All types).
# These two sheets list overlapping events, so both receive updates.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 ---
$ws1.Range("F2").Value = 193
$ws1.Range("F3").Value = 5402
$ws1.Range("F5").Value = 55
$ws1.Range("F7").Value = 617
$ws1.Range("F8").Value = 590
$ws1.Range("F9").Value = 1057
$ws1.Range("F11").Value = 1483
$ws1.Range("F12").Value = 4425
$ws1.Range("F14").Value = 197
$ws1.Range("F15").Value = 171
$ws1.Range("F17").Value = 3496
$ws1.Range("F18").Value = 176
$ws1.Range("F19").Value = 1108
$ws1.Range("F20").Value = 107
$ws1.Range("F22").Value = 202
$ws1.Range("F23").Value = 19
$ws1.Range("F24").Value = 130
$ws1.Range("F25").Value = 44
$ws1.Range("F26").Value = 143
$ws1.Range("F29").Value = 32
$ws1.Range("F30").Value = 57
$ws1.Range("F32").Value = 30
$ws1.Range("F33").Value = 31

# --- Sheet 4: 全部类型 ---
$ws4.Range("F2").Value = 193
$ws4.Range("F4").Value = 5403
$ws4.Range("F6").Value = 55
$ws4.Range("F8").Value = 617
$ws4.Range("F9").Value = 590
$ws4.Range("F10").Value = 1057
$ws4.Range("F12").Value = 1483
$ws4.Range("F13").Value = 4425
$ws4.Range("F15").Value = 197
$ws4.Range("F16").Value = 171
$ws4.Range("F18").Value = 3496
$ws4.Range("F19").Value = 176
$ws4.Range("F20").Value = 1108
$ws4.Range("F21").Value = 107
$ws4.Range("F23").Value = 202
$ws4.Range("F24").Value = 19
$ws4.Range("F25").Value = 130
$ws4.Range("F26").Value = 44
$ws4.Range("F27").Value = 143
$ws4.Range("F30").Value = 32
$ws4.Range("F31").Value = 57
$ws4.Range("F33").Value = 30
$ws4.Range("F34").Value = 31

Write-Host "Updated counts applied."
